$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.464.79"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "1.909.98"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'238.83"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.4773"
$ws.Range("E7").Value = "  -2.30%  "
$ws.Range("D8").Value = "'0.2829"
$ws.Range("E8").Value = "  -3.60%  "
$ws.Range("D9").Value = "'0.06698"
$ws.Range("E9").Value = "  -3.38%  "
$ws.Range("D10").Value = "'18.65"
$ws.Range("E10").Value = "  -4.41%  "
$ws.Range("D11").Value = "'100.68"
$ws.Range("E11").Value = "  -4.83%  "
$ws.Range("D12").Value = "1.917.63"
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("D13").Value = "'0.07673"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").Value = "'5.195"
$ws.Range("D15").Value = "'0.6667"
$ws.Range("E15").Value = "  -4.69%  "
$ws.Range("D16").Value = "30.492.81"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "'255.70"
$ws.Range("E17").Value = "  -7.12%  "
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "'0.000007448"
$ws.Range("E19").Value = "  -3.88%  "
$ws.Range("D20").Value = "'12.62"
$ws.Range("E20").Value = "  -3.80%  "
$ws.Range("D21").Value = "'5.373"
$ws.Range("E21").Value = "  -1.18%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "'6.280"
$ws.Range("E23").Value = "  -3.28%  "
$ws.Range("D24").Value = "'9.326"
$ws.Range("E24").Value = "  -4.07%  "
$ws.Range("D25").Value = "'167.66"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("D26").Value = "'19.06"
$ws.Range("E26").Value = "  -2.78%  "
$ws.Range("D27").Value = "'2.052"
$ws.Range("E27").Value = "  -4.83%  "
$ws.Range("D28").Value = "'4.766"
$ws.Range("E28").Value = "  +4.60%  "
$ws.Range("D29").Value = "'1.386"
$ws.Range("D30").Value = "'0.1001"
$ws.Range("E30").Value = "  -3.76%  "
$ws.Range("D31").Value = "'1.507"
$ws.Range("E31").Value = "  -2.90%  "
$ws.Range("D32").Value = "'4.247"
$ws.Range("E32").Value = "  -2.70%  "
$ws.Range("D33").Value = "'0.04700"
$ws.Range("E33").Value = "  -3.32%  "
$ws.Range("D34").Value = "'0.7230"
$ws.Range("E34").Value = "  -3.62%  "
$ws.Range("D35").Value = "'1.103"
$ws.Range("E35").Value = "  -4.53%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.701"
$ws.Range("E36").Value = "  -0.88%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.01910"
$ws.Range("E37").Value = "  -4.15%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.610"
$ws.Range("E38").Value = "  -2.01%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'6.262"
$ws.Range("E39").Value = "  -3.19%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "'74.97"
$ws.Range("E40").Value = "  -3.72%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'1.959"
$ws.Range("E41").Value = "  -6.84%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.8600"
$ws.Range("E42").Value = "  -4.39%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'105.35"
$ws.Range("E43").Value = "  -2.90%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'0.9998"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.4229"
$ws.Range("E45").Value = "  -3.94%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "'7.359"
$ws.Range("E46").Value = "  -4.72%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.1196"
$ws.Range("E47").Value = "  -3.85%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "'918.67"
$ws.Range("E48").Value = "  -6.98%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'34.67"
$ws.Range("E49").Value = "  -3.05%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.754"
$ws.Range("E50").Value = "  -5.59%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05738"
$ws.Range("E51").Value = "  +0.17%  "
